$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 3, shifting existing rows 3-4 down to 4-5.
$ws.Rows.Item(3).Insert()

# The newly inserted row 3 duplicates row 2 (the "iaest-measure:..." row).
$ws.Range("A3").Value = $ws.Range("A2").Formula
$ws.Range("B3").Value = $ws.Range("B2").Formula
$ws.Range("C3").Value = $ws.Range("C2").Formula
$ws.Range("D3").Value = $ws.Range("D2").Formula
$ws.Range("E3").Value = $ws.Range("E2").Formula
